# Generate Report for Handoff
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps for the rows whose handoff xliff was just (re)generated.
# - Mark those same rows' Priority as "ht" (handoff type) on both the
#   zh-cn and de-de target sheets.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 10, 11, 13, 14)

# Overview sheet: column G is "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-12 18:25:18"
}

# zh-cn sheet: column H is "Latest Handoff Datetime", column E is "Priority"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-12 18:25:00"
    $wsZhCn.Range("E$r").Value = "ht"
}

# de-de sheet: column H is "Latest Handoff Datetime", column E is "Priority"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-12 18:25:18"
    $wsDeDe.Range("E$r").Value = "ht"
}
